$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = 14.14129999999999
$ws.Range("E10").Value = 12.3306
$ws.Range("E12").Value = 12.367
$ws.Range("E18").Value = 12.5087
